# completed data extraction and data merge
#
# The "Sheet1" tab lists each NC county's population estimates; column A held
# values like "Alamance County". Strip the redundant " County" suffix so the
# column just shows the county name (e.g. "Alamance"), matching the cleaned
# data used elsewhere in the merge. Also leave the selection on the sheet
# parked on column N (a scratch/next-free column) as it was when the edit was
# saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows for the 100 counties live in A2:A101 (row 1 is the "County"
# header, row 102 is the North Carolina statewide total).
$firstRow = 2
$lastRow  = 101
$rowCount = $lastRow - $firstRow + 1

$dataRange = $ws.Range("A$firstRow" + ":A$lastRow")

# Read the current text for every county row first...
$names = New-Object "object[,]" $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $cell = $ws.Cells.Item($firstRow + $i, 1)
    $text = $cell.Text
    if ($text.EndsWith(" County")) {
        $text = $text.Substring(0, $text.Length - 7)
    }
    $names[$i, 0] = $text
}

# ...then write all the cleaned names back in a single batched assignment so
# the shared-string table is rebuilt once, in original order, instead of
# shifting around on every individual cell write.
$dataRange.Value2 = $names

# Restore the sheet's last-known selection: the whole of column N.
[void]$ws.Activate()
[void]$ws.Range("N:N").Select()
